$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared text strings (volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/18/2025  Through  8/24/2025"

# --- Cells that switch from a numeric value to the text sentinel "0" (style 13, shared string) ---
# Use Range.Copy() from a cell that already carries the target text-style formatting,
# which preserves both style (s=13) and shared-string reuse (value "0").
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("C28"))

# --- Cells that switch from the text sentinel to a real number (style 14, integer-like format) ---
$ws.Range("I14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("I14").Copy($ws.Range("G27"))
$ws.Range("G27").Value = 1
$ws.Range("I14").Copy($ws.Range("G28"))
$ws.Range("G28").Value = 1

# --- Cells that switch from the text sentinel to a real number (style 15, percent-change format) ---
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("K14").Copy($ws.Range("H27"))
$ws.Range("H27").Value = 0
$ws.Range("K14").Copy($ws.Range("H28"))
$ws.Range("H28").Value = 100

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -97.368421052631
$ws.Range("M15").Value = -61.904761904761
$ws.Range("N15").Value = -77.142857142857
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -61.538461538461
$ws.Range("I16").Value = 61
$ws.Range("J16").Value = 86
$ws.Range("K16").Value = -29.069767441860
$ws.Range("L16").Value = -21.794871794871
$ws.Range("M16").Value = -57.342657342657
$ws.Range("N16").Value = -86.563876651982
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -47.619047619047
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 153
$ws.Range("K17").Value = -31.372549019607
$ws.Range("L17").Value = -31.818181818181
$ws.Range("M17").Value = -16.666666666666
$ws.Range("N17").Value = -78.170478170478
$ws.Range("D18").Value = 7
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -83.333333333333
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = 4.761904761904
$ws.Range("M18").Value = -30.158730158730
$ws.Range("N18").Value = -91.002044989775
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -19.230769230769
$ws.Range("I19").Value = 234
$ws.Range("J19").Value = 226
$ws.Range("K19").Value = 3.539823008849
$ws.Range("M19").Value = 160
$ws.Range("N19").Value = 9.345794392523
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 48
$ws.Range("J20").Value = 47
$ws.Range("K20").Value = 2.127659574468
$ws.Range("L20").Value = -34.246575342465
$ws.Range("M20").Value = 41.176470588235
$ws.Range("N20").Value = -76.119402985074
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -75
$ws.Range("F21").Value = 46
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = -42.5
$ws.Range("I21").Value = 501
$ws.Range("J21").Value = 587
$ws.Range("K21").Value = -14.650766609880
$ws.Range("L21").Value = -12.105263157894
$ws.Range("M21").Value = 4.158004158004
$ws.Range("N21").Value = -73.797071129707
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -25
$ws.Range("M22").Value = 114.285714285714
$ws.Range("M23").Value = -60
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").Value = 46
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = -33.333333333333
$ws.Range("I24").Value = 500
$ws.Range("J24").Value = 497
$ws.Range("K24").Value = 0.603621730382
$ws.Range("L24").Value = -5.482041587901
$ws.Range("M24").Value = 116.450216450216
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -81.818181818181
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = -79.310344827586
$ws.Range("I25").Value = 95
$ws.Range("J25").Value = 142
$ws.Range("K25").Value = -33.098591549295
$ws.Range("L25").Value = -31.159420289855
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = -12.195121951219
$ws.Range("I26").Value = 237
$ws.Range("J26").Value = 234
$ws.Range("K26").Value = 1.282051282051
$ws.Range("L26").Value = 15.609756097561
$ws.Range("M26").Value = -26.168224299065
$ws.Range("I27").Value = 12
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = 9.090909090909
$ws.Range("L27").Value = 9.090909090909
$ws.Range("F28").Value = 2
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = 40
$ws.Range("L28").Value = -4.545454545454
$ws.Range("M29").Value = -83.333333333333
$ws.Range("N29").Value = -96.610169491525
$ws.Range("M30").Value = -80
$ws.Range("N30").Value = -96.491228070175
